# Applies the "massive reworking of the UI" edit to joint.xlsx.
# Sheet tab "Sheet2" is the big 50-row joint table (xl/worksheets/sheet1.xml);
# Sheet tab "Sheet1" is the small 6-row bone-chain table (xl/worksheets/sheet2.xml).

$wb = $excel.ActiveWorkbook
$big = $wb.Worksheets.Item("Sheet2")
$small = $wb.Worksheets.Item("Sheet1")

$xlNone = -4142
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Big sheet: rows 15-32 used to carry a style (s="7", borderId=0 applyBorder)
#    on several B cells and on a block of now-empty D:I / P:U cells. Strip the
#    border from those B cells (collapses them back to the default style) and
#    outright clear the dead D:I / P:U cells so they disappear from the XML.
# ---------------------------------------------------------------------------
$bCellsToUnstyle = @("B15","B16","B21","B22","B27","B28")
foreach ($addr in $bCellsToUnstyle) {
    $big.Range($addr).Borders.LineStyle = $xlNone
}

$big.Range("D15:I32").Clear()
$big.Range("P15:U32").Clear()

# ---------------------------------------------------------------------------
# 2. Big sheet: rows 33-50 gain an explicit C column value of 0, and the J
#    column loses its stray border-only style (s="7" -> default).
# ---------------------------------------------------------------------------
for ($r = 33; $r -le 50; $r++) {
    $big.Cells.Item($r, 3).Value = 0
    $big.Range("J" + $r).Borders.LineStyle = $xlNone
}

# ---------------------------------------------------------------------------
# 3. Style-table cleanup: cellXfs 7 (border-only, applyFont duplicate of 4)
#    and cellXfs 8 (border-only variant of 4) / 9 (fill+border variant of 5)
#    are redundant with 4 and 5. Re-point every cell that used the redundant
#    "applyFont" styles (8/9) at the equivalent canonical style (4/5) by
#    copy/pasting formats from a cell that already carries the canonical
#    style, so Excel collapses the now-unused xf entries out of styles.xml.
# ---------------------------------------------------------------------------
$style4Donor = $big.Range("A17")   # plain border style (cellXfs index 4)
$style5Donor = $big.Range("A18")   # fill + border style (cellXfs index 5)

$toStyle4 = @("A39","B42","B43","B44","A43","A44","A45","B48","B49","A49","A50","B50")
$toStyle5 = @("B39","A40","B40","A41","B41","A42","B45","A46","B46","A47","B47","A48")

foreach ($addr in $toStyle4) {
    $style4Donor.Copy() | Out-Null
    $big.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}
foreach ($addr in $toStyle5) {
    $style5Donor.Copy() | Out-Null
    $big.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

$smallToStyle4 = @("A1","B4","A5","B5","A6","B6")
$smallToStyle5 = @("B1","A2","B2","A3","B3","A4")

foreach ($addr in $smallToStyle4) {
    $style4Donor.Copy() | Out-Null
    $small.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}
foreach ($addr in $smallToStyle5) {
    $style5Donor.Copy() | Out-Null
    $small.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. View state: scroll position + active selection moved during the UI pass.
# ---------------------------------------------------------------------------
$bigWindow = $excel.ActiveWindow
$bigWindow.ScrollRow = 17
$big.Range("C49").Select()

Write-Host "done"
